$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Block 1: refresh rows 2-23 (sector/stock summary rows) on "Recommandations" ---
$data1 = New-Object "object[,]" 22,7
$data1[0,0] = 'BRVM - SERVICES PUBLICS'
$data1[0,1] = 0.0
$data1[0,2] = 8.0
$data1[0,3] = 3253.85
$data1[0,4] = 105.01
$data1[0,5] = '🟡 Observer'
$data1[0,6] = '➖ Neutre'
$data1[1,0] = 'SAFCA CI'
$data1[1,1] = 0.0
$data1[1,2] = 4.0
$data1[1,3] = 3095.0
$data1[1,4] = 860.0
$data1[1,5] = '🟡 Observer'
$data1[1,6] = '➖ Neutre'
$data1[2,0] = 'BRVM - AUTRES SECTEURS'
$data1[2,1] = 0.0
$data1[2,2] = 4.0
$data1[2,3] = 2614.76
$data1[2,4] = 658.16
$data1[2,5] = '🟡 Observer'
$data1[2,6] = '➖ Neutre'
$data1[3,0] = 'CFAO MOTORS CI'
$data1[3,1] = 0.0
$data1[3,2] = 4.0
$data1[3,3] = 2610.0
$data1[3,4] = 655.0
$data1[3,5] = '🟡 Observer'
$data1[3,6] = '➖ Neutre'
$data1[4,0] = 'SETAO CI'
$data1[4,1] = 0.0
$data1[4,2] = 4.0
$data1[4,3] = 2325.0
$data1[4,4] = 590.0
$data1[4,5] = '🟡 Observer'
$data1[4,6] = '➖ Neutre'
$data1[5,0] = 'NEI-CEDA CI'
$data1[5,1] = 0.0
$data1[5,2] = 4.0
$data1[5,3] = 2300.0
$data1[5,4] = 595.0
$data1[5,5] = '🟡 Observer'
$data1[5,6] = '➖ Neutre'
$data1[6,0] = 'UNIWAX CI'
$data1[6,1] = 0.0
$data1[6,2] = 4.0
$data1[6,3] = 2255.0
$data1[6,4] = 575.0
$data1[6,5] = '🟡 Observer'
$data1[6,6] = '➖ Neutre'
$data1[7,0] = 'AIR LIQUIDE CI'
$data1[7,1] = 0.0
$data1[7,2] = 4.0
$data1[7,3] = 2060.0
$data1[7,4] = 525.0
$data1[7,5] = '🟡 Observer'
$data1[7,6] = '➖ Neutre'
$data1[8,0] = 'BRVM - DISTRIBUTION'
$data1[8,1] = 0.0
$data1[8,2] = 4.0
$data1[8,3] = 1438.22
$data1[8,4] = 362.92
$data1[8,5] = '🟡 Observer'
$data1[8,6] = '➖ Neutre'
$data1[9,0] = 'BRVM - TRANSPORT'
$data1[9,1] = 0.0
$data1[9,2] = 4.0
$data1[9,3] = 1393.97
$data1[9,4] = 352.47
$data1[9,5] = '🟡 Observer'
$data1[9,6] = '➖ Neutre'
$data1[10,0] = 'BRVM - AGRICULTURE'
$data1[10,1] = 0.0
$data1[10,2] = 4.0
$data1[10,3] = 1255.17
$data1[10,4] = 316.23
$data1[10,5] = '🟡 Observer'
$data1[10,6] = '➖ Neutre'
$data1[11,0] = 'BRVM - INDUSTRIE'
$data1[11,1] = 0.0
$data1[11,2] = 4.0
$data1[11,3] = 1013.84
$data1[11,4] = 261.41
$data1[11,5] = '🟡 Observer'
$data1[11,6] = '➖ Neutre'
$data1[12,0] = 'BRVM - CONSOMMATION DE BASE'
$data1[12,1] = 0.0
$data1[12,2] = 4.0
$data1[12,3] = 837.77
$data1[12,4] = 215.21
$data1[12,5] = '🟡 Observer'
$data1[12,6] = '➖ Neutre'
$data1[13,0] = 'BRVM-PRINCIPAL'
$data1[13,1] = 0.0
$data1[13,2] = 4.0
$data1[13,3] = 747.65
$data1[13,4] = 189.01
$data1[13,5] = '🟡 Observer'
$data1[13,6] = '➖ Neutre'
$data1[14,0] = 'BRVM - INDUSTRIELS'
$data1[14,1] = 0.0
$data1[14,2] = 4.0
$data1[14,3] = 549.52
$data1[14,4] = 138.9
$data1[14,5] = '🟡 Observer'
$data1[14,6] = '➖ Neutre'
$data1[15,0] = 'BRVM-PRESTIGE'
$data1[15,1] = 0.0
$data1[15,2] = 4.0
$data1[15,3] = 517.73
$data1[15,4] = 129.27
$data1[15,5] = '🟡 Observer'
$data1[15,6] = '➖ Neutre'
$data1[16,0] = 'BRVM - FINANCES'
$data1[16,1] = 0.0
$data1[16,2] = 4.0
$data1[16,3] = 488.87
$data1[16,4] = 121.93
$data1[16,5] = '🟡 Observer'
$data1[16,6] = '➖ Neutre'
$data1[17,0] = 'BRVM - SERVICES FINANCIERS'
$data1[17,1] = 0.0
$data1[17,2] = 4.0
$data1[17,3] = 480.46
$data1[17,4] = 119.83
$data1[17,5] = '🟡 Observer'
$data1[17,6] = '➖ Neutre'
$data1[18,0] = 'BRVM - CONSOMMATION DISCRETIONNAIRE'
$data1[18,1] = 0.0
$data1[18,2] = 4.0
$data1[18,3] = 424.48
$data1[18,4] = 105.5
$data1[18,5] = '🟡 Observer'
$data1[18,6] = '➖ Neutre'
$data1[19,0] = 'BRVM - ENERGIE'
$data1[19,1] = 0.0
$data1[19,2] = 4.0
$data1[19,3] = 421.49
$data1[19,4] = 107.45
$data1[19,5] = '🟡 Observer'
$data1[19,6] = '➖ Neutre'
$data1[20,0] = 'BRVM - TELECOMMUNICATIONS'
$data1[20,1] = 0.0
$data1[20,2] = 4.0
$data1[20,3] = 369.27
$data1[20,4] = 92.57
$data1[20,5] = '🟡 Observer'
$data1[20,6] = '➖ Neutre'
$data1[21,0] = 'UNILEVER CI (UNLC)'
$data1[21,1] = 4.0
$data1[21,2] = 0.0
$data1[21,3] = 29.87
$data1[21,4] = 7.49
$data1[21,5] = '🟢 Achat'
$data1[21,6] = '✅ Renforcer'
$ws1.Range("A2:G23").Value = $data1

# --- Block 2: replace the stock ranking table (rows 28-49), sorted desc by "Variation Totale" ---
# Original rows 28-51 covered 24 stocks; SOGB CI (SOGC) and SONATEL SN (SNTS) are dropped
# and the remaining 22 stocks are re-sorted/refreshed, shrinking the table to rows 28-49.
$data2 = New-Object "object[,]" 22,7
$data2[0,0] = 'SOCIETE GENERALE COTE D''IVOIRE (SGBC)'
$data2[0,1] = 1
$data2[0,2] = 1
$data2[0,3] = 3.58
$data2[0,4] = 7.04
$data2[0,5] = '🟡 Observer'
$data2[0,6] = '👀 À surveiller'
$data2[1,0] = 'SETAO CI (STAC)'
$data2[1,1] = 1
$data2[1,2] = 0
$data2[1,3] = 3.48
$data2[1,4] = 3.48
$data2[1,5] = '🟡 Observer'
$data2[1,6] = '➖ Neutre'
$data2[2,0] = 'SOCIETE IVOIRIENNE DE BANQUE  (SIBC)'
$data2[2,1] = 2
$data2[2,2] = 1
$data2[2,3] = 2.88
$data2[2,4] = -3.03
$data2[2,5] = '🟡 Observer'
$data2[2,6] = '👀 À surveiller'
$data2[3,0] = 'BANK OF AFRICA BN (BOAB)'
$data2[3,1] = 1
$data2[3,2] = 0
$data2[3,3] = 2.86
$data2[3,4] = 2.86
$data2[3,5] = '🟡 Observer'
$data2[3,6] = '➖ Neutre'
$data2[4,0] = 'BANK OF AFRICA SENEGAL (BOAS)'
$data2[4,1] = 1
$data2[4,2] = 1
$data2[4,3] = 2.43
$data2[4,4] = -2.22
$data2[4,5] = '🟡 Observer'
$data2[4,6] = '👀 À surveiller'
$data2[5,0] = 'FILTISAC CI (FTSC)'
$data2[5,1] = 1
$data2[5,2] = 1
$data2[5,3] = 2.09
$data2[5,4] = -0.8
$data2[5,5] = '🟡 Observer'
$data2[5,6] = '👀 À surveiller'
$data2[6,0] = 'CIE CI (CIEC)'
$data2[6,1] = 1
$data2[6,2] = 1
$data2[6,3] = 0.51
$data2[6,4] = -5.06
$data2[6,5] = '🟡 Observer'
$data2[6,6] = '👀 À surveiller'
$data2[7,0] = 'TOTAL'
$data2[7,1] = 0
$data2[7,2] = 4
$data2[7,3] = 0
$data2[7,4] = 0
$data2[7,5] = '🟡 Observer'
$data2[7,6] = '➖ Neutre'
$data2[8,0] = 'UNIWAX CI (UNXC)'
$data2[8,1] = 1
$data2[8,2] = 1
$data2[8,3] = -0.71
$data2[8,4] = 3.6
$data2[8,5] = '🟡 Observer'
$data2[8,6] = '👀 À surveiller'
$data2[9,0] = 'SOLIBRA CI (SLBC)'
$data2[9,1] = 0
$data2[9,2] = 1
$data2[9,3] = -0.81
$data2[9,4] = -0.81
$data2[9,5] = '🟡 Observer'
$data2[9,6] = '➖ Neutre'
$data2[10,0] = 'NEI-CEDA CI (NEIC)'
$data2[10,1] = 0
$data2[10,2] = 1
$data2[10,3] = -0.84
$data2[10,4] = -0.84
$data2[10,5] = '🟡 Observer'
$data2[10,6] = '➖ Neutre'
$data2[11,0] = 'TOTALENERGIES MARKETING CI (TTLC)'
$data2[11,1] = 0
$data2[11,2] = 1
$data2[11,3] = -1.6
$data2[11,4] = -1.6
$data2[11,5] = '🟡 Observer'
$data2[11,6] = '➖ Neutre'
$data2[12,0] = 'ORAGROUP TOGO (ORGT)'
$data2[12,1] = 0
$data2[12,2] = 1
$data2[12,3] = -1.74
$data2[12,4] = -1.74
$data2[12,5] = '🟡 Observer'
$data2[12,6] = '➖ Neutre'
$data2[13,0] = 'VIVO ENERGY CI (SHEC)'
$data2[13,1] = 0
$data2[13,2] = 1
$data2[13,3] = -2
$data2[13,4] = -2
$data2[13,5] = '🟡 Observer'
$data2[13,6] = '➖ Neutre'
$data2[14,0] = 'ONATEL BF (ONTBF)'
$data2[14,1] = 0
$data2[14,2] = 1
$data2[14,3] = -2.08
$data2[14,4] = -2.08
$data2[14,5] = '🟡 Observer'
$data2[14,6] = '➖ Neutre'
$data2[15,0] = 'SITAB CI (STBC)'
$data2[15,1] = 0
$data2[15,2] = 1
$data2[15,3] = -2.5
$data2[15,4] = -2.5
$data2[15,5] = '🟡 Observer'
$data2[15,6] = '➖ Neutre'
$data2[16,0] = 'SICOR CI (SICC)'
$data2[16,1] = 0
$data2[16,2] = 1
$data2[16,3] = -2.73
$data2[16,4] = -2.73
$data2[16,5] = '🟡 Observer'
$data2[16,6] = '➖ Neutre'
$data2[17,0] = 'TRACTAFRIC MOTORS CI (PRSC)'
$data2[17,1] = 0
$data2[17,2] = 1
$data2[17,3] = -3.42
$data2[17,4] = -3.42
$data2[17,5] = '🟡 Observer'
$data2[17,6] = '➖ Neutre'
$data2[18,0] = 'BERNABE CI (BNBC)'
$data2[18,1] = 1
$data2[18,2] = 1
$data2[18,3] = -3.53
$data2[18,4] = 3.59
$data2[18,5] = '🟡 Observer'
$data2[18,6] = '👀 À surveiller'
$data2[19,0] = 'CFAO MOTORS CI (CFAC)'
$data2[19,1] = 0
$data2[19,2] = 1
$data2[19,3] = -3.65
$data2[19,4] = -3.65
$data2[19,5] = '🟡 Observer'
$data2[19,6] = '➖ Neutre'
$data2[20,0] = 'SUCRIVOIRE (SCRC)'
$data2[20,1] = 0
$data2[20,2] = 1
$data2[20,3] = -3.85
$data2[20,4] = -3.85
$data2[20,5] = '🟡 Observer'
$data2[20,6] = '➖ Neutre'
$data2[21,0] = 'BANK OF AFRICA NG (BOAN)'
$data2[21,1] = 0
$data2[21,2] = 2
$data2[21,3] = -4.19
$data2[21,4] = -2.19
$data2[21,5] = '🟡 Observer'
$data2[21,6] = '➖ Neutre'
$ws1.Range("A28:G49").Value = $data2

# --- Remove the two now-obsolete trailing rows (SOGB CI (SOGC) / SONATEL SN (SNTS) dropped) ---
$ws1.Rows.Item(51).Delete()
$ws1.Rows.Item(50).Delete()

# --- Block 3: refresh "Top_YTD" volume column (B2:B11) ---
$data3 = New-Object "object[,]" 10,1
$data3[0,0] = 7546780.15
$data3[1,0] = 576663.2
$data3[2,0] = 322539.33
$data3[3,0] = 320369.19
$data3[4,0] = 215262.8
$data3[5,0] = 207185.49
$data3[6,0] = 193882.34
$data3[7,0] = 142878.52
$data3[8,0] = 44499.47
$data3[9,0] = 40352.47
$ws2.Range("B2:B11").Value = $data3
